# Quarterly financials update for ALEX sheet.
# Two new fiscal quarters (period ending 2018-12-31 and 2018-09-30) are
# inserted as new columns D and E; the pre-existing quarterly columns
# (old D:K) shift right to F:M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALEX")

# Insert two new blank columns at D, shifting the old D:K data to F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# The freshly inserted D:E columns come back with the formatting of the
# column that used to be immediately to their left (General / column C).
# Re-apply the real column formatting (date format row, #,##0 data rows,
# right alignment, etc.) by copying it over from the columns that hold
# the equivalent, now-shifted, historical quarters (F:G).
$ws.Range("F5:G102").Copy() | Out-Null
$ws.Range("D5:E102").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# New quarter-end values for every row that carries data, as
# @(row, D-value, E-value).
$newData = @(
    @(7, 43465, 43373),
    @(8, 299600, 119400),
    @(9, 114900, 87100),
    @(10, 184700, 32300),
    @(12, "NA", "NA"),
    @(13, 0, 0),
    @(14, 79400, "NA"),
    @(15, 0, 0),
    @(17, 159400, 101700),
    @(18, 140200, 17700),
    @(20, -198800, 8200),
    @(21, -47400, 36200),
    @(22, 8900, 9100),
    @(23, -67500, 16800),
    @(24, 18100, 1000),
    @(25, 0, 0),
    @(26, -85600, 15800),
    @(27, -136200, 15000),
    @(28, 0, 0),
    @(29, -400, -200),
    @(30, 0, 0),
    @(31, 0, 0),
    @(32, 198800, -8200),
    @(33, -136600, 14800),
    @(34, 0, 0),
    @(35, -136600, 14800),
    @(38, 43465, 43373),
    @(41, 11400, 7500),
    @(42, 0, 0),
    @(43, 75000, 60300),
    @(44, 26500, 32200),
    @(45, 67800, 74600),
    @(46, 180700, 174600),
    @(47, 171400, 379200),
    @(48, 1322000, 1322300),
    @(49, 133500, 178200),
    @(50, 0, 0),
    @(51, 0, 0),
    @(52, 417600, 221900),
    @(53, 0, 0),
    @(54, 2225200, 2276200),
    @(57, 34200, 33700),
    @(58, 39000, 37100),
    @(59, 56200, 49900),
    @(60, 129400, 120700),
    @(61, 739100, 741300),
    @(62, 140500, 57600),
    @(63, 0, 0),
    @(64, 0, 0),
    @(65, 0, 0),
    @(66, 1022600, 932900),
    @(68, 0, 0),
    @(69, 0, 0),
    @(70, 0, 0),
    @(71, 0, 0),
    @(72, -538900, -410500),
    @(73, 0, 0),
    @(74, 0, 0),
    @(75, 0, 0),
    @(76, 1202600, 1343300),
    @(77, 0, 0),
    @(80, 43465, 43373),
    @(81, -136600, 14800),
    @(83, 11200, 10300),
    @(84, 0, 0),
    @(85, 0, 0),
    @(86, 0, 0),
    @(87, 0, 0),
    @(88, 0, 0),
    @(89, 272200, 10900),
    @(91, -14400, -14700),
    @(92, 0, 0),
    @(93, 0, 0),
    @(94, -43900, -600),
    @(96, 0, 0),
    @(97, 0, 0),
    @(98, 0, 0),
    @(99, 0, 0),
    @(100, -1100, -12300),
    @(101, 0, 0),
    @(102, 227200, -2000)
)

foreach ($entry in $newData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 4).Value = $entry[1]   # column D
    $ws.Cells.Item($r, 5).Value = $entry[2]   # column E
}

# Row 91 ("Capital Expenditures") also carries a historical restatement
# for three of the older quarters beyond the simple column shift.
$ws.Range("G91").Value = -12700
$ws.Range("I91").Value = -23600
$ws.Range("J91").Value = -4000
